$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Profile_pic_URL" header (column P) to "Profile picture"
$ws.Range("P1").Value = "Profile picture"

# Add the new "Remark" field as column Q
$ws.Range("Q1").Value = "Remark"

# Move the active selection to the newly added header cell
$ws.Range("Q1").Select()
